# Auto-generated script applying the Tonberry_Profits market-data refresh diff.
# Updates currentAveragePrice* / Leve*Price / Leve*Profit columns (H:N) per row
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 989.2727
$ws.Range("I40").Value = 1052.875
$ws.Range("J40").Value = 819.6667
$ws.Range("K40").Value = 1052.875
$ws.Range("L40").Value = 819.6667
$ws.Range("M40").Value = -877.875
$ws.Range("N40").Value = -1169.6667
$ws.Range("H106").Value = 1665.875
$ws.Range("I106").Value = 1665.875
$ws.Range("K106").Value = 1665.875
$ws.Range("M106").Value = -1034.875
$ws.Range("H118").Value = 529.6667
$ws.Range("I118").Value = 529.6667
$ws.Range("K118").Value = 1589.0001
$ws.Range("M118").Value = 67.99990000000003
$ws.Range("H129").Value = 877.4318
$ws.Range("J129").Value = 892.95
$ws.Range("L129").Value = 2678.85
$ws.Range("N129").Value = -12678.85
$ws.Range("H135").Value = 373.5926
$ws.Range("I135").Value = 385.65384
$ws.Range("J135").Value = 60
$ws.Range("K135").Value = 3470.88456
$ws.Range("L135").Value = 540
$ws.Range("M135").Value = -935.88456
$ws.Range("N135").Value = -5610
$ws.Range("H137").Value = 1347.2759
$ws.Range("I137").Value = 836.2917
$ws.Range("J137").Value = 3800
$ws.Range("K137").Value = 2508.8751
$ws.Range("L137").Value = 11400
$ws.Range("M137").Value = 41.1248999999998
$ws.Range("N137").Value = -16500
$ws.Range("H138").Value = 4735.1304
$ws.Range("I138").Value = 10459.4
$ws.Range("J138").Value = 3145.0557
$ws.Range("K138").Value = 31378.2
$ws.Range("L138").Value = 9435.167099999999
$ws.Range("M138").Value = -26238.2
$ws.Range("N138").Value = -19715.1671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3392.6875
$ws.Range("I32").Value = 3005.35
$ws.Range("K32").Value = 3005.35
$ws.Range("M32").Value = -2718.35
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H45").Value = 1344.1875
$ws.Range("I45").Value = 1139.238
$ws.Range("J45").Value = 1735.4546
$ws.Range("K45").Value = 1139.238
$ws.Range("L45").Value = 1735.4546
$ws.Range("M45").Value = -762.2380000000001
$ws.Range("N45").Value = -2489.4546
$ws.Range("H61").Value = 2389.2163
$ws.Range("I61").Value = 1527.303
$ws.Range("J61").Value = 9500
$ws.Range("K61").Value = 1527.303
$ws.Range("L61").Value = 9500
$ws.Range("M61").Value = -1315.303
$ws.Range("N61").Value = -9924
$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 1000
$ws.Range("K102").Value = 1000
$ws.Range("M102").Value = 622
$ws.Range("H132").Value = 1336.1842
$ws.Range("I132").Value = 943.0294
$ws.Range("J132").Value = 4678
$ws.Range("K132").Value = 2829.0882
$ws.Range("L132").Value = 14034
$ws.Range("M132").Value = -299.0882000000001
$ws.Range("N132").Value = -19094
$ws.Range("H136").Value = 2389.2163
$ws.Range("I136").Value = 1527.303
$ws.Range("J136").Value = 9500
$ws.Range("K136").Value = 4581.909000000001
$ws.Range("L136").Value = 28500
$ws.Range("M136").Value = -2031.909000000001
$ws.Range("N136").Value = -33600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 185871.64
$ws.Range("I86").Value = 7649.75
$ws.Range("J86").Value = 287712.72
$ws.Range("K86").Value = 7649.75
$ws.Range("L86").Value = 287712.72
$ws.Range("M86").Value = -6526.75
$ws.Range("N86").Value = -289958.72
$ws.Range("H89").Value = 185871.64
$ws.Range("I89").Value = 7649.75
$ws.Range("J89").Value = 287712.72
$ws.Range("K89").Value = 38248.75
$ws.Range("L89").Value = 1438563.6
$ws.Range("M89").Value = -32632.75
$ws.Range("N89").Value = -1449795.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1644.0322
$ws.Range("I31").Value = 1502.4584
$ws.Range("J31").Value = 2129.4285
$ws.Range("K31").Value = 1502.4584
$ws.Range("L31").Value = 2129.4285
$ws.Range("M31").Value = -1207.4584
$ws.Range("N31").Value = -2719.4285
$ws.Range("H34").Value = 1644.0322
$ws.Range("I34").Value = 1502.4584
$ws.Range("J34").Value = 2129.4285
$ws.Range("K34").Value = 1502.4584
$ws.Range("L34").Value = 2129.4285
$ws.Range("M34").Value = -1300.4584
$ws.Range("N34").Value = -2533.4285
$ws.Range("H58").Value = 1740535.1
$ws.Range("I58").Value = 2174998.2
$ws.Range("K58").Value = 2174998.2
$ws.Range("M58").Value = -2174795.2
$ws.Range("H99").Value = 2402.3
$ws.Range("I99").Value = 2431.8572
$ws.Range("K99").Value = 2431.8572
$ws.Range("M99").Value = -933.8571999999999
$ws.Range("H122").Value = 2853.6875
$ws.Range("J122").Value = 3916.875
$ws.Range("L122").Value = 11750.625
$ws.Range("N122").Value = -16650.625
$ws.Range("H126").Value = 2402.3
$ws.Range("I126").Value = 2431.8572
$ws.Range("K126").Value = 7295.571599999999
$ws.Range("M126").Value = -4825.571599999999
$ws.Range("H132").Value = 1436.909
$ws.Range("I132").Value = 880.3929000000001
$ws.Range("K132").Value = 2641.1787
$ws.Range("M132").Value = -111.1787000000004
$ws.Range("H134").Value = 1225.3414
$ws.Range("I134").Value = 1147.7646
$ws.Range("J134").Value = 1602.1428
$ws.Range("K134").Value = 3443.2938
$ws.Range("L134").Value = 4806.428400000001
$ws.Range("M134").Value = -908.2937999999999
$ws.Range("N134").Value = -9876.428400000001
$ws.Range("H136").Value = 1740535.1
$ws.Range("I136").Value = 2174998.2
$ws.Range("K136").Value = 6524994.600000001
$ws.Range("M136").Value = -6522444.600000001
$ws.Range("H141").Value = 73332
$ws.Range("J141").Value = 73332
$ws.Range("L141").Value = 73332
$ws.Range("N141").Value = -83692

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3274.875
$ws.Range("J55").Value = 3274.875
$ws.Range("L55").Value = 9824.625
$ws.Range("N55").Value = -10178.625
$ws.Range("H75").Value = 1006.5
$ws.Range("J75").Value = 2000
$ws.Range("L75").Value = 6000
$ws.Range("N75").Value = -7996
$ws.Range("H78").Value = 1006.5
$ws.Range("J78").Value = 2000
$ws.Range("L78").Value = 18000
$ws.Range("N78").Value = -27984
$ws.Range("H132").Value = 943.1429000000001
$ws.Range("J132").Value = 1149.5
$ws.Range("L132").Value = 10345.5
$ws.Range("N132").Value = -15405.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4385411.5
$ws.Range("J7").Value = 1010399.6
$ws.Range("L7").Value = 1010399.6
$ws.Range("N7").Value = -1010623.6
$ws.Range("H8").Value = 4385411.5
$ws.Range("J8").Value = 1010399.6
$ws.Range("L8").Value = 1010399.6
$ws.Range("N8").Value = -1010677.6
$ws.Range("H14").Value = 2840000.8
$ws.Range("I14").Value = 3425001
$ws.Range("J14").Value = 500000
$ws.Range("K14").Value = 3425001
$ws.Range("L14").Value = 500000
$ws.Range("M14").Value = -3424833
$ws.Range("N14").Value = -500336
$ws.Range("H102").Value = 2538.9678
$ws.Range("I102").Value = 2579.24
$ws.Range("J102").Value = 2371.1667
$ws.Range("K102").Value = 2579.24
$ws.Range("L102").Value = 2371.1667
$ws.Range("M102").Value = -957.2399999999998
$ws.Range("N102").Value = -5615.1667
$ws.Range("H122").Value = 1769.6154
$ws.Range("I122").Value = 1581
$ws.Range("J122").Value = 1931.2858
$ws.Range("K122").Value = 4743
$ws.Range("L122").Value = 5793.857400000001
$ws.Range("M122").Value = -2293
$ws.Range("N122").Value = -10693.8574
$ws.Range("H132").Value = 1167346
$ws.Range("I132").Value = 1673552
$ws.Range("J132").Value = 3072.4
$ws.Range("K132").Value = 5020656
$ws.Range("L132").Value = 9217.200000000001
$ws.Range("M132").Value = -5018126
$ws.Range("N132").Value = -14277.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3349.4
$ws.Range("I22").Value = 5997.5
$ws.Range("J22").Value = 2687.375
$ws.Range("K22").Value = 5997.5
$ws.Range("L22").Value = 2687.375
$ws.Range("M22").Value = -5702.5
$ws.Range("N22").Value = -3277.375
$ws.Range("H27").Value = 3349.4
$ws.Range("I27").Value = 5997.5
$ws.Range("J27").Value = 2687.375
$ws.Range("K27").Value = 5997.5
$ws.Range("L27").Value = 2687.375
$ws.Range("M27").Value = -5890.5
$ws.Range("N27").Value = -2901.375
$ws.Range("H100").Value = 1563.3334
$ws.Range("I100").Value = 1625
$ws.Range("J100").Value = 1440
$ws.Range("K100").Value = 1625
$ws.Range("L100").Value = 1440
$ws.Range("M100").Value = -1084
$ws.Range("N100").Value = -2522
$ws.Range("H122").Value = 3122.2
$ws.Range("I122").Value = 1985.6666
$ws.Range("J122").Value = 7668.3335
$ws.Range("K122").Value = 5956.9998
$ws.Range("L122").Value = 23005.0005
$ws.Range("M122").Value = -3506.9998
$ws.Range("N122").Value = -27905.0005
$ws.Range("H132").Value = 1570.1154
$ws.Range("I132").Value = 1070
$ws.Range("J132").Value = 2153.5833
$ws.Range("K132").Value = 3210
$ws.Range("L132").Value = 6460.749899999999
$ws.Range("M132").Value = -680
$ws.Range("N132").Value = -11520.7499
$ws.Range("H136").Value = 2120.3333
$ws.Range("I136").Value = 1226.4546
$ws.Range("K136").Value = 3679.3638
$ws.Range("M136").Value = -1129.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 46788.41
$ws.Range("I122").Value = 49637.688
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 148913.064
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -146463.064
$ws.Range("N122").Value = -8500
$ws.Range("H124").Value = 23746.5
$ws.Range("J124").Value = 23746.5
$ws.Range("L124").Value = 23746.5
$ws.Range("N124").Value = -33566.5
$ws.Range("H126").Value = 5304.75
$ws.Range("I126").Value = 2781.6
$ws.Range("J126").Value = 7107
$ws.Range("K126").Value = 8344.799999999999
$ws.Range("L126").Value = 21321
$ws.Range("M126").Value = -5874.799999999999
$ws.Range("N126").Value = -26261
$ws.Range("H132").Value = 2408.3333
$ws.Range("I132").Value = 1356.8462
$ws.Range("K132").Value = 4070.5386
$ws.Range("M132").Value = -1540.5386
$ws.Range("H136").Value = 17923182
$ws.Range("I136").Value = 25254308
$ws.Range("J136").Value = 2656.111
$ws.Range("K136").Value = 75762924
$ws.Range("L136").Value = 7968.333
$ws.Range("M136").Value = -75760374
$ws.Range("N136").Value = -13068.333

Write-Output "Applied 263 cell updates across 8 sheets."
